# Add the missing trial-distribution data for the last participant (row 12,
# "VH11") across the five data sheets, and update each sheet's last-saved
# selection / active-tab bookkeeping to match.
#
# Column map for row 12 (J:M are newly populated; A:I already had data):
#   nTrialsCollected : 10 10 10 10
#   nNoResponses     :  1  0  0  1
#   nHandError       :  0  1  1  0
#   nEquipmentErrors :  0  0  0  0
#   nTrialsAnalysed  :  9  9  9  9

$wb = $excel.ActiveWorkbook

# --- nTrialsCollected ---------------------------------------------------
$ws = $wb.Worksheets.Item("nTrialsCollected")
$ws.Range("J12").Value = 10
$ws.Range("K12").Value = 10
$ws.Range("L12").Value = 10
$ws.Range("M12").Value = 10
$ws.Range("N12").Select() | Out-Null

# --- nNoResponses ---------------------------------------------------------
$ws = $wb.Worksheets.Item("nNoResponses")
$ws.Range("J12").Value = 1
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 1
$ws.Range("L13").Select() | Out-Null

# --- nHandError -------------------------------------------------------
$ws = $wb.Worksheets.Item("nHandError")
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0
$ws.Range("K13").Select() | Out-Null

# --- nTrialsAnalysed --------------------------------------------------
$ws = $wb.Worksheets.Item("nTrialsAnalysed")
$ws.Range("J12").Value = 9
$ws.Range("K12").Value = 9
$ws.Range("L12").Value = 9
$ws.Range("M12").Value = 9
$ws.Range("M13").Select() | Out-Null

# --- nEquipmentErrors ---------------------------------------------------
# This is the sheet left active/selected when the workbook was last saved
# (activeTab moves from nTrialsCollected to nEquipmentErrors), so it must
# be the last sheet activated in this script.
$ws = $wb.Worksheets.Item("nEquipmentErrors")
$ws.Activate()
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("M13").Select() | Out-Null
